# Pre_Site_Total_Alarms.xlsx update
# - Add a "Date" header label in A1
# - Reorder the site columns C:G from (AZC, PKV, SDU, STL, THL)
#   to (PKV, STL, THL, AZC, SDU)
# - Data rows 2-5 follow the same column reorder (values per site unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row values, columns A through G
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "ARN"
$ws.Range("C1").Value = "PKV"
$ws.Range("D1").Value = "STL"
$ws.Range("E1").Value = "THL"
$ws.Range("F1").Value = "AZC"
$ws.Range("G1").Value = "SDU"

# Row 2
$ws.Range("C2").Value = 330
$ws.Range("D2").Value = 44
$ws.Range("E2").Value = 45
$ws.Range("F2").Value = 580
$ws.Range("G2").Value = 73

# Row 3
$ws.Range("C3").Value = 345
$ws.Range("D3").Value = 45
$ws.Range("E3").Value = 92
$ws.Range("F3").Value = 625
$ws.Range("G3").Value = 84

# Row 4
$ws.Range("C4").Value = 345
$ws.Range("D4").Value = 45
$ws.Range("E4").Value = 152
$ws.Range("F4").Value = 674
$ws.Range("G4").Value = 86

# Row 5
$ws.Range("C5").Value = 459
$ws.Range("D5").Value = 208
$ws.Range("E5").Value = 269
$ws.Range("F5").Value = 756
$ws.Range("G5").Value = 93
